$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames (accents stripped)
$ws.Range("B1").Value = "Frequencia"
$ws.Range("C1").Value = "Matematica"
$ws.Range("D1").Value = "Ciencias"
$ws.Range("E1").Value = "Fisica"

# Row 2 - Vinicius
$ws.Range("B2").Value = 260
$ws.Range("D2").Value = 5.5

# Row 3 - Marco
$ws.Range("B3").Value = 245
$ws.Range("C3").Value = 8.5

# Row 4 - Vivian
$ws.Range("B4").Value = 255
$ws.Range("E4").Value = 7.5

# Row 5 - Ana
$ws.Range("B5").Value = 230
$ws.Range("C5").Value = 8.5
$ws.Range("D5").Value = 7.5

# Row 6 - Tatiana
$ws.Range("B6").Value = 259
